$wb = $excel.ActiveWorkbook

# --- Rename "Sheet1" to "category" ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "category"

# --- Fill olist_category / master_category mapping data ---
$data = New-Object 'object[,]' 72,2
$data[0,0] = 'olist_category'
$data[0,1] = 'master_category'
$data[1,0] = 'bed_bath_table'
$data[1,1] = 'Meuble Déco'
$data[2,0] = 'sports_leisure'
$data[2,1] = 'Sport'
$data[3,0] = 'health_beauty'
$data[3,1] = 'Mode Bagage'
$data[4,0] = 'furniture_decor'
$data[4,1] = 'Meuble Déco'
$data[5,0] = 'computers_accessories'
$data[5,1] = 'Informatique'
$data[6,0] = 'toys'
$data[6,1] = 'Jouets'
$data[7,0] = 'housewares'
$data[7,1] = 'Electoménager'
$data[8,0] = 'cool_stuff'
$data[8,1] = 'Mode Bagage'
$data[9,0] = 'watches_gifts'
$data[9,1] = 'Mode Bagage'
$data[10,0] = 'telephony'
$data[10,1] = 'Téléphonie'
$data[11,0] = 'garden_tools'
$data[11,1] = 'Jardin Animalerie'
$data[12,0] = 'perfumery'
$data[12,1] = 'Mode Bagage'
$data[13,0] = 'auto'
$data[13,1] = 'Auto Moto'
$data[14,0] = 'baby'
$data[14,1] = 'Bébé'
$data[15,0] = 'fashion_bags_accessories'
$data[15,1] = 'Mode Bagage'
$data[16,0] = 'stationery'
$data[16,1] = 'Librairie'
$data[17,0] = 'electronics'
$data[17,1] = 'Informatique'
$data[18,0] = 'pet_shop'
$data[18,1] = 'Jardin Animalerie'
$data[19,0] = 'consoles_games'
$data[19,1] = 'Jeux Vidéo'
$data[20,0] = 'office_furniture'
$data[20,1] = 'Meuble Déco'
$data[21,0] = 'luggage_accessories'
$data[21,1] = 'Mode Bagage'
$data[22,0] = 'small_appliances'
$data[22,1] = 'Electoménager'
$data[23,0] = 'musical_instruments'
$data[23,1] = 'TV Son Photo'
$data[24,0] = 'home_confort'
$data[24,1] = 'Meuble Déco'
$data[25,0] = 'home_appliances'
$data[25,1] = 'Electoménager'
$data[26,0] = 'books_general_interest'
$data[26,1] = 'Librairie'
$data[27,0] = 'furniture_living_room'
$data[27,1] = 'Meuble Déco'
$data[28,0] = 'market_place'
$data[28,1] = 'Autres'
$data[29,0] = 'fashion_shoes'
$data[29,1] = 'Mode Bagage'
$data[30,0] = 'audio'
$data[30,1] = 'TV Son Photo'
$data[31,0] = 'computers'
$data[31,1] = 'Informatique'
$data[32,0] = 'fixed_telephony'
$data[32,1] = 'Téléphonie'
$data[33,0] = 'home_construction'
$data[33,1] = 'Bricolage'
$data[34,0] = 'air_conditioning'
$data[34,1] = 'Bricolage'
$data[35,0] = 'kitchen_dining_laundry_garden_furniture'
$data[35,1] = 'Jardin Animalerie'
$data[36,0] = 'food_drink'
$data[36,1] = 'Alimentation Boisson'
$data[37,0] = 'construction_tools_construction'
$data[37,1] = 'Bricolage'
$data[38,0] = 'food'
$data[38,1] = 'Alimentation Boisson'
$data[39,0] = 'home_appliances_2'
$data[39,1] = 'Electoménager'
$data[40,0] = 'fashion_male_clothing'
$data[40,1] = 'Mode Bagage'
$data[41,0] = 'fashion_underwear_beach'
$data[41,1] = 'Mode Bagage'
$data[42,0] = 'drinks'
$data[42,1] = 'Alimentation Boisson'
$data[43,0] = 'costruction_tools_garden'
$data[43,1] = 'Jardin Animalerie'
$data[44,0] = 'tablets_printing_image'
$data[44,1] = 'Informatique'
$data[45,0] = 'christmas_supplies'
$data[45,1] = 'Meuble Déco'
$data[46,0] = 'agro_industry_and_commerce'
$data[46,1] = 'Autres'
$data[47,0] = 'books_technical'
$data[47,1] = 'Librairie'
$data[48,0] = 'furniture_bedroom'
$data[48,1] = 'Meuble Déco'
$data[49,0] = 'dvds_blu_ray'
$data[49,1] = 'TV Son Photo'
$data[50,0] = 'construction_tools_safety'
$data[50,1] = 'Bricolage'
$data[51,0] = 'signaling_and_security'
$data[51,1] = 'Bricolage'
$data[52,0] = 'industry_commerce_and_business'
$data[52,1] = 'Autres'
$data[53,0] = 'art'
$data[53,1] = 'Meuble Déco'
$data[54,0] = 'fashio_female_clothing'
$data[54,1] = 'Mode Bagage'
$data[55,0] = 'costruction_tools_tools'
$data[55,1] = 'Bricolage'
$data[56,0] = 'fashion_sport'
$data[56,1] = 'Sport'
$data[57,0] = 'furniture_mattress_and_upholstery'
$data[57,1] = 'Meuble Déco'
$data[58,0] = 'home_comfort_2'
$data[58,1] = 'Meuble Déco'
$data[59,0] = 'construction_tools_lights'
$data[59,1] = 'Bricolage'
$data[60,0] = 'books_imported'
$data[60,1] = 'Librairie'
$data[61,0] = 'cds_dvds_musicals'
$data[61,1] = 'TV Son Photo'
$data[62,0] = 'music'
$data[62,1] = 'TV Son Photo'
$data[63,0] = 'party_supplies'
$data[63,1] = 'Meuble Déco'
$data[64,0] = 'la_cuisine'
$data[64,1] = 'Electoménager'
$data[65,0] = 'cine_photo'
$data[65,1] = 'TV Son Photo'
$data[66,0] = 'flowers'
$data[66,1] = 'Jardin Animalerie'
$data[67,0] = 'fashion_childrens_clothes'
$data[67,1] = 'Mode Bagage'
$data[68,0] = 'small_appliances_home_oven_and_coffee'
$data[68,1] = 'Electoménager'
$data[69,0] = 'security_and_services'
$data[69,1] = 'Informatique'
$data[70,0] = 'diapers_and_hygiene'
$data[70,1] = 'Bébé'
$data[71,0] = 'arts_and_craftmanship'
$data[71,1] = 'Meuble Déco'
$ws2.Range("A1:B72").Value = $data

# --- Column widths (bestFit) ---
$ws2.Columns.Item(1).AutoFit() | Out-Null
$ws2.Columns.Item(2).AutoFit() | Out-Null

# --- Turn the range into an Excel Table ---
$lo = $ws2.ListObjects.Add(1, $ws2.Range("A1:B72"), 0, 1)
$lo.Name = "Table2"
$lo.TableStyle = "TableStyleMedium2"

# --- Page setup: landscape, paper size 9 (A4) ---
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 2

# --- Make "category" the active / selected sheet (matches activeTab + tabSelected) ---
$ws2.Activate() | Out-Null
$ws2.Range("D7").Select() | Out-Null
